$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.056.12"
$ws.Range("E2").Value = "  +3.71%  "
$ws.Range("D3").Value = "1.893.75"
$ws.Range("E3").Value = "  +3.99%  "
$ws.Range("D4").Value = "'0.9979"
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").Value = "'247.50"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").Value = "'0.9981"
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("D7").Value = "'0.4983"
$ws.Range("E7").Value = "  +1.15%  "
$ws.Range("D8").Value = "'44.82"
$ws.Range("E8").Value = "  +1.23%  "
$ws.Range("D9").Value = "'0.2954"
$ws.Range("E9").Value = "  +6.58%  "
$ws.Range("D10").Value = "'0.06657"
$ws.Range("E10").Value = "  +4.27%  "
$ws.Range("D11").Value = "1.893.16"
$ws.Range("E11").Value = "  +3.85%  "
$ws.Range("D12").Value = "'17.03"
$ws.Range("E12").Value = "  +1.92%  "
$ws.Range("D13").Value = "'0.07227"
$ws.Range("E13").Value = "  +2.24%  "
$ws.Range("D14").Value = "'0.6800"
$ws.Range("E14").Value = "  +5.81%  "
$ws.Range("D15").Value = "'85.89"
$ws.Range("E15").Value = "  +2.16%  "
$ws.Range("D16").Value = "'4.870"
$ws.Range("E16").Value = "  +3.97%  "
$ws.Range("D17").Value = "30.027.03"
$ws.Range("E17").Value = "  +3.50%  "
$ws.Range("D18").Value = "'0.000007994"
$ws.Range("E18").Value = "  +9.45%  "
$ws.Range("D19").Value = "'0.9984"
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("D20").Value = "'12.97"
$ws.Range("E20").Value = "  +6.12%  "
$ws.Range("D21").Value = "2.136.40"
$ws.Range("E21").Value = "  +4.13%  "
$ws.Range("D22").Value = "'0.9972"
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("D23").Value = "'4.784"
$ws.Range("E23").Value = "  +5.21%  "
$ws.Range("D24").Value = "'5.693"
$ws.Range("E24").Value = "  +6.11%  "
$ws.Range("D25").Value = "'9.220"
$ws.Range("E25").Value = "  +4.48%  "
$ws.Range("D26").Value = "'147.40"
$ws.Range("E26").Value = "  +2.21%  "
$ws.Range("D27").Value = "'131.91"
$ws.Range("E27").Value = "  +2.75%  "
$ws.Range("D28").Value = "'16.84"
$ws.Range("E28").Value = "  +2.78%  "
$ws.Range("D29").Value = "'1.953"
$ws.Range("E29").Value = "  +3.84%  "
$ws.Range("D30").Value = "'1.364"
$ws.Range("E30").Value = "  -2.81%  "
$ws.Range("D31").Value = "'4.260"
$ws.Range("E31").Value = "  +3.38%  "
$ws.Range("D32").Value = "'0.08766"
$ws.Range("E32").Value = "  +4.88%  "
$ws.Range("D33").Value = "'3.963"
$ws.Range("E33").Value = "  +5.03%  "
$ws.Range("D34").Value = "'0.05097"
$ws.Range("E34").Value = "  +2.97%  "
$ws.Range("D35").Value = "'1.120"
$ws.Range("E35").Value = "  +2.17%  "
$ws.Range("D36").Value = "'0.7071"
$ws.Range("E36").Value = "  +5.39%  "
$ws.Range("D37").Value = "'2.666"
$ws.Range("E37").Value = "  -0.97%  "
$ws.Range("D38").Value = "'2.780"
$ws.Range("E38").Value = "  +3.79%  "
$ws.Range("D39").Value = "'2.220"
$ws.Range("E39").Value = "  -2.77%  "
$ws.Range("D40").Value = "'0.9443"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").Value = "'0.01665"
$ws.Range("E41").Value = "  +5.16%  "
$ws.Range("D42").Value = "'6.070"
$ws.Range("E42").Value = "  -2.11%  "
$ws.Range("D43").Value = "'0.9968"
$ws.Range("E43").Value = "  -0.38%  "
$ws.Range("D44").Value = "'103.49"
$ws.Range("E44").Value = "  +2.50%  "
$ws.Range("D45").Value = "'0.4220"
$ws.Range("E45").Value = "  +3.80%  "
$ws.Range("D46").Value = "'7.504"
$ws.Range("E46").Value = "  +5.00%  "
$ws.Range("D47").Value = "'0.1262"
$ws.Range("E47").Value = "  +3.47%  "
$ws.Range("D48").Value = "'0.05723"
$ws.Range("E48").Value = "  +3.48%  "
$ws.Range("D49").Value = "'32.83"
$ws.Range("E49").Value = "  +3.65%  "
$ws.Range("D50").Value = "'8.278"
$ws.Range("E50").Value = "  +2.37%  "
$ws.Range("D51").Value = "'0.3743"
$ws.Range("E51").Value = "  +3.96%  "
